# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / Leve profit figures
# from a scheduled data-refresh run across several Leve sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 222.63158
$ws.Cells.Item(33, 9).Value = 180
$ws.Cells.Item(33, 10).Value = 295.7143
$ws.Cells.Item(33, 11).Value = 180
$ws.Cells.Item(33, 12).Value = 295.7143
$ws.Cells.Item(33, 13).Value = 49
$ws.Cells.Item(33, 14).Value = -753.7143

# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 18126026
$ws.Cells.Item(40, 9).Value = 62500000
$ws.Cells.Item(40, 10).Value = 3334700.8
$ws.Cells.Item(40, 11).Value = 62500000
$ws.Cells.Item(40, 12).Value = 3334700.8
$ws.Cells.Item(40, 13).Value = -62499825
$ws.Cells.Item(40, 14).Value = -3335050.8

# Row 76 (Leve Item ID 12602)
$ws.Cells.Item(76, 8).Value = 9809906
$ws.Cells.Item(76, 9).Value = 7103.96
$ws.Cells.Item(76, 10).Value = 37039910
$ws.Cells.Item(76, 11).Value = 7103.96
$ws.Cells.Item(76, 12).Value = 37039910
$ws.Cells.Item(76, 13).Value = -6788.96
$ws.Cells.Item(76, 14).Value = -37040540

# Row 79 (Leve Item ID 12602)
$ws.Cells.Item(79, 8).Value = 9809906
$ws.Cells.Item(79, 9).Value = 7103.96
$ws.Cells.Item(79, 10).Value = 37039910
$ws.Cells.Item(79, 11).Value = 7103.96
$ws.Cells.Item(79, 12).Value = 37039910
$ws.Cells.Item(79, 13).Value = -6011.96
$ws.Cells.Item(79, 14).Value = -37042094

# Row 111 (Leve Item ID 27768)
$ws.Cells.Item(111, 8).Value = 636.6667
$ws.Cells.Item(111, 9).Value = 449.66666
$ws.Cells.Item(111, 10).Value = 1010.6667
$ws.Cells.Item(111, 11).Value = 1348.99998
$ws.Cells.Item(111, 12).Value = 3032.0001
$ws.Cells.Item(111, 13).Value = 1718.00002
$ws.Cells.Item(111, 14).Value = -9166.000100000001

# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 7250389
$ws.Cells.Item(132, 9).Value = 8199841.5
$ws.Cells.Item(132, 10).Value = 10809.875
$ws.Cells.Item(132, 11).Value = 24599524.5
$ws.Cells.Item(132, 12).Value = 32429.625
$ws.Cells.Item(132, 13).Value = -24596994.5
$ws.Cells.Item(132, 14).Value = -37489.625

$ws = $wb.Worksheets.Item("ARM")
# Row 8 (Leve Item ID 3011)
$ws.Cells.Item(8, 8).Value = 6002300
$ws.Cells.Item(8, 9).Value = 10000500
$ws.Cells.Item(8, 10).Value = 5000
$ws.Cells.Item(8, 11).Value = 10000500
$ws.Cells.Item(8, 12).Value = 5000
$ws.Cells.Item(8, 13).Value = -10000356
$ws.Cells.Item(8, 14).Value = -5288

# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 11631419
$ws.Cells.Item(32, 9).Value = 3593.25
$ws.Cells.Item(32, 10).Value = 55558760
$ws.Cells.Item(32, 11).Value = 3593.25
$ws.Cells.Item(32, 12).Value = 55558760
$ws.Cells.Item(32, 13).Value = -3306.25
$ws.Cells.Item(32, 14).Value = -55559334

# Row 35 (Leve Item ID 2473)
$ws.Cells.Item(35, 8).Value = 1468.5
$ws.Cells.Item(35, 9).Value = 1468.5
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1468.5
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1062.5
$ws.Cells.Item(35, 14).ClearContents()

# Row 61 (Leve Item ID 43999)
$ws.Cells.Item(61, 8).Value = 1858.5
$ws.Cells.Item(61, 9).Value = 1858.5
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1858.5
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1646.5
$ws.Cells.Item(61, 14).ClearContents()

# Row 136 (Leve Item ID 43999)
$ws.Cells.Item(136, 8).Value = 1858.5
$ws.Cells.Item(136, 9).Value = 1858.5
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 5575.5
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -3025.5
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Cells.Item(22, 8).Value = 13513513
$ws.Cells.Item(22, 9).Value = 13513513
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 13513513
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -13513340

# Row 23 (Leve Item ID 1686)
$ws.Cells.Item(23, 8).Value = 5000
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 5000
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 5000
$ws.Cells.Item(23, 14).Value = -5566

# Row 29 (Leve Item ID 2318)
$ws.Cells.Item(29, 8).Value = 1007.5
$ws.Cells.Item(29, 9).Value = 1007.5
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 1007.5
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -718.5
$ws.Cells.Item(29, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 25 (Leve Item ID 4709)
$ws.Cells.Item(25, 8).Value = 300
$ws.Cells.Item(25, 9).Value = 300
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 900
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -731
$ws.Cells.Item(25, 14).ClearContents()

# Row 26 (Leve Item ID 4746)
$ws.Cells.Item(26, 8).Value = 30340
$ws.Cells.Item(26, 9).Value = 100
$ws.Cells.Item(26, 10).Value = 50500
$ws.Cells.Item(26, 11).Value = 300
$ws.Cells.Item(26, 12).Value = 151500
$ws.Cells.Item(26, 13).Value = -12
$ws.Cells.Item(26, 14).Value = -152076

# Row 30 (Leve Item ID 4709)
$ws.Cells.Item(30, 8).Value = 300
$ws.Cells.Item(30, 9).Value = 300
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 900
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = -798
$ws.Cells.Item(30, 14).ClearContents()

# Row 69 (Leve Item ID 12850)
$ws.Cells.Item(69, 8).Value = 7924.1577
$ws.Cells.Item(69, 9).Value = 500
$ws.Cells.Item(69, 10).Value = 8797.588
$ws.Cells.Item(69, 11).Value = 1500
$ws.Cells.Item(69, 12).Value = 26392.764
$ws.Cells.Item(69, 13).Value = -689
$ws.Cells.Item(69, 14).Value = -28014.764

# Row 72 (Leve Item ID 12850)
$ws.Cells.Item(72, 8).Value = 7924.1577
$ws.Cells.Item(72, 9).Value = 500
$ws.Cells.Item(72, 10).Value = 8797.588
$ws.Cells.Item(72, 11).Value = 4500
$ws.Cells.Item(72, 12).Value = 79178.292
$ws.Cells.Item(72, 13).Value = -444
$ws.Cells.Item(72, 14).Value = -87290.292

# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 1229.27
$ws.Cells.Item(131, 9).Value = 606
$ws.Cells.Item(131, 10).Value = 1262.0737
$ws.Cells.Item(131, 11).Value = 1818
$ws.Cells.Item(131, 12).Value = 3786.2211
$ws.Cells.Item(131, 13).Value = 3222
$ws.Cells.Item(131, 14).Value = -13866.2211

$ws = $wb.Worksheets.Item("GSM")
# Row 3 (Leve Item ID 4091)
$ws.Cells.Item(3, 8).Value = 617.6667
$ws.Cells.Item(3, 9).Value = 402
$ws.Cells.Item(3, 10).Value = 833.3333
$ws.Cells.Item(3, 11).Value = 402
$ws.Cells.Item(3, 12).Value = 833.3333
$ws.Cells.Item(3, 13).Value = -286
$ws.Cells.Item(3, 14).Value = -1065.3333

# Row 5 (Leve Item ID 1681)
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).ClearContents()

# Row 9 (Leve Item ID 1683)
$ws.Cells.Item(9, 8).Value = 5000
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 5000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 14).Value = -5340

# Row 12 (Leve Item ID 4093)
$ws.Cells.Item(12, 8).Value = 3002565.8
$ws.Cells.Item(12, 9).Value = 4201592
$ws.Cells.Item(12, 10).Value = 5000
$ws.Cells.Item(12, 11).Value = 4201592
$ws.Cells.Item(12, 12).Value = 5000
$ws.Cells.Item(12, 13).Value = -4201452
$ws.Cells.Item(12, 14).Value = -5280

# Row 14 (Leve Item ID 4198)
$ws.Cells.Item(14, 8).Value = 787.375
$ws.Cells.Item(14, 9).Value = 185.57143
$ws.Cells.Item(14, 10).Value = 5000
$ws.Cells.Item(14, 11).Value = 185.57143
$ws.Cells.Item(14, 12).Value = 5000
$ws.Cells.Item(14, 13).Value = -17.57142999999999
$ws.Cells.Item(14, 14).Value = -5336

$ws = $wb.Worksheets.Item("LTW")
# Row 19 (Leve Item ID 2229)
$ws.Cells.Item(19, 8).Value = 200
$ws.Cells.Item(19, 9).Value = 200
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 200
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -30

# Row 21 (Leve Item ID 2672)
$ws.Cells.Item(21, 8).Value = 4835.6665
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 4835.6665
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 4835.6665
$ws.Cells.Item(21, 14).Value = -5183.6665
